# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (want-to-go count) values in the F column on both
# the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 813
$ws1.Range("F8").Value = 4639
$ws1.Range("F10").Value = 5064
$ws1.Range("F11").Value = 579
$ws1.Range("F12").Value = 1270

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 813
$ws4.Range("F9").Value = 4639
$ws4.Range("F11").Value = 5064
$ws4.Range("F12").Value = 579
$ws4.Range("F13").Value = 1270
